$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order of assignment controls the shared-string table order, so cells are
# written in the exact sequence needed to reproduce the target sharedStrings
# index layout (0..12).

# Header row (row 1) - rename "from" component labels -> indices 0..3
$ws.Range("B1").Value = "Q_from_net1"
$ws.Range("C1").Value = "Q_from_CHP1"
$ws.Range("D1").Value = "Q_from_solar_th1"
$ws.Range("E1").Value = "Q_from_pvt1"

# "to" row labels -> indices 4..5
$ws.Range("A2").Value = "Q_to_demand1"
$ws.Range("A3").Value = "Q_to_net1"

# Row 2 data (to demand1) -> indices 6..8 (then 9 for D2)
$ws.Range("B2").Value = "Q_net1_demand1"
$ws.Range("C2").Value = "Q_CHP1_demand1"
$ws.Range("C3").Value = "Q_CHP1_net1"
$ws.Range("D2").Value = "Q_solar_th1_demand1"
$ws.Range("D3").Value = "Q_solar_th1_net1"
$ws.Range("E2").Value = "Q_pvt1_demand1"
$ws.Range("E3").Value = "Q_pvt1_net1"

# Row 3, column B stays a plain numeric zero (unchanged)
$ws.Range("B3").Value = 0
